# Update "Fase de Grupos" (Group Stage) sheet with the results of the
# remaining Group D and Group F matches, and correct the kickoff dates
# of the Round-of-16 matches (Jogo 55 and Jogo 57).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fase de Grupos")

# The sheet is protected (data-entry cells are unlocked, but the date
# cells we also need to touch are locked) - unprotect first so all the
# writes succeed, then restore protection afterwards.
$ws.Unprotect("CC01")

# --- Grupo D ---------------------------------------------------------
# Los Angeles FC 1 x 1 Flamengo
$ws.Range("F25").Value2 = 1
$ws.Range("H25").Value2 = 1

# Esperance de Tunis 0 x 3 Chelsea
$ws.Range("F26").Value2 = 0
$ws.Range("H26").Value2 = 3

# --- Grupo F ---------------------------------------------------------
# Mamelodi Sundowns 0 x 0 Fluminense
$ws.Range("F37").Value2 = 0
$ws.Range("H37").Value2 = 0

# Borussia Dortmund 1 x 0 Ulsan
$ws.Range("F38").Value2 = 1
$ws.Range("H38").Value2 = 0

# --- Datas dos jogos das oitavas (Round of 16) ------------------------
# Jogo 55
$ws.Range("C43").Value2 = 45834
$ws.Range("C44").Value2 = 45834
# Jogo 57
$ws.Range("C49").Value2 = 45834
$ws.Range("C50").Value2 = 45834

# Restore sheet protection as it was before the edit.
$ws.Protect("CC01", $true, $true, $true)

# --- View state --------------------------------------------------------
# Reflect where the user was looking/selecting when they saved.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F31").Select()

$wb.Save()
